$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the now-obsolete TODO list item descriptions (B3, B4, B5, B7) while
# keeping "Overhaul commandHandler" in B6 and the title row intact.
$ws.Range("B3").ClearContents()
$ws.Range("B4").ClearContents()
$ws.Range("B5").ClearContents()
$ws.Range("B7").ClearContents()

# Rows 3 and 7 had custom heights to accommodate the long wrapped text that
# was just removed; let Excel recompute the row heights now that the
# content is gone.
$ws.Rows("3").AutoFit()
$ws.Rows("7").AutoFit()
